# Auto-generated edit script: update Leve market-price columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()

$ws.Cells.Item(58, 8).Value = 1892.7858
$ws.Cells.Item(58, 9).Value = 254.57143
$ws.Cells.Item(58, 10).Value = 3531
$ws.Cells.Item(58, 11).Value = 763.71429
$ws.Cells.Item(58, 12).Value = 10593
$ws.Cells.Item(58, 13).Value = -613.71429
$ws.Cells.Item(58, 14).Value = -10893

$ws.Cells.Item(103, 8).Value = 320
$ws.Cells.Item(103, 9).Value = 320
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 960
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -374
$ws.Cells.Item(103, 14).ClearContents()

$ws.Cells.Item(116, 8).Value = 3647.5
$ws.Cells.Item(116, 9).Value = 3647.5
$ws.Cells.Item(116, 11).Value = 3647.5
$ws.Cells.Item(116, 13).Value = -205.5

$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()

$ws.Cells.Item(138, 8).Value = 1893.8334
$ws.Cells.Item(138, 10).Value = 2737.125
$ws.Cells.Item(138, 12).Value = 8211.375
$ws.Cells.Item(138, 14).Value = -18491.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1917.125
$ws.Cells.Item(32, 9).Value = 1202
$ws.Cells.Item(32, 11).Value = 1202
$ws.Cells.Item(32, 13).Value = -915

$ws.Cells.Item(37, 8).Value = 9375

$ws.Cells.Item(40, 8).Value = 25171.334
$ws.Cells.Item(40, 9).Value = 1028
$ws.Cells.Item(40, 10).Value = 30000
$ws.Cells.Item(40, 11).Value = 1028
$ws.Cells.Item(40, 12).Value = 30000
$ws.Cells.Item(40, 13).Value = -852
$ws.Cells.Item(40, 14).Value = -30352

$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()

$ws.Cells.Item(55, 8).Value = 30000
$ws.Cells.Item(55, 10).Value = 30000
$ws.Cells.Item(55, 12).Value = 30000
$ws.Cells.Item(55, 14).Value = -30630

$ws.Cells.Item(80, 8).Value = 23333.334

$ws.Cells.Item(83, 8).Value = 23333.334

$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()

$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3178.1428
$ws.Cells.Item(99, 9).Value = 2649.6
$ws.Cells.Item(99, 10).Value = 4499.5
$ws.Cells.Item(99, 11).Value = 2649.6
$ws.Cells.Item(99, 12).Value = 4499.5
$ws.Cells.Item(99, 13).Value = -1151.6
$ws.Cells.Item(99, 14).Value = -7495.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1058.4
$ws.Cells.Item(16, 9).Value = 1058.4
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1058.4
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -771.4000000000001
$ws.Cells.Item(16, 14).ClearContents()

$ws.Cells.Item(19, 8).Value = 4521.5
$ws.Cells.Item(19, 9).Value = 216.25
$ws.Cells.Item(19, 10).Value = 8826.75
$ws.Cells.Item(19, 11).Value = 216.25
$ws.Cells.Item(19, 12).Value = 8826.75
$ws.Cells.Item(19, 13).Value = -46.25
$ws.Cells.Item(19, 14).Value = -9166.75

$ws.Cells.Item(24, 8).Value = 4521.5
$ws.Cells.Item(24, 9).Value = 216.25
$ws.Cells.Item(24, 10).Value = 8826.75
$ws.Cells.Item(24, 11).Value = 216.25
$ws.Cells.Item(24, 12).Value = 8826.75
$ws.Cells.Item(24, 13).Value = -46.25
$ws.Cells.Item(24, 14).Value = -9166.75

$ws.Cells.Item(47, 8).Value = 86537.5
$ws.Cells.Item(47, 9).Value = 12000
$ws.Cells.Item(47, 10).Value = 111383.336
$ws.Cells.Item(47, 11).Value = 12000
$ws.Cells.Item(47, 12).Value = 111383.336
$ws.Cells.Item(47, 13).Value = -11434
$ws.Cells.Item(47, 14).Value = -112515.336

$ws.Cells.Item(113, 8).Value = 1058.4
$ws.Cells.Item(113, 9).Value = 1058.4
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1058.4
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 1111.6
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 2813.647
$ws.Cells.Item(134, 10).Value = 4459.5
$ws.Cells.Item(134, 12).Value = 13378.5
$ws.Cells.Item(134, 14).Value = -18448.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 2417
$ws.Cells.Item(17, 9).Value = 1800
$ws.Cells.Item(17, 10).Value = 2681.4285
$ws.Cells.Item(17, 11).Value = 5400
$ws.Cells.Item(17, 12).Value = 8044.2855
$ws.Cells.Item(17, 13).Value = -5231
$ws.Cells.Item(17, 14).Value = -8382.2855

$ws.Cells.Item(26, 8).Value = 50.75
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 12).Value = 150
$ws.Cells.Item(26, 14).Value = -726

$ws.Cells.Item(34, 8).Value = 902.4
$ws.Cells.Item(34, 10).Value = 1345.6666
$ws.Cells.Item(34, 12).Value = 4036.9998
$ws.Cells.Item(34, 14).Value = -4204.9998

$ws.Cells.Item(41, 8).Value = 236
$ws.Cells.Item(41, 10).Value = 236
$ws.Cells.Item(41, 12).Value = 708
$ws.Cells.Item(41, 14).Value = -1384

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 2965.5715
$ws.Cells.Item(22, 9).Value = 50
$ws.Cells.Item(22, 10).Value = 3451.5
$ws.Cells.Item(22, 11).Value = 50
$ws.Cells.Item(22, 12).Value = 3451.5
$ws.Cells.Item(22, 13).Value = 479
$ws.Cells.Item(22, 14).Value = -4509.5

$ws.Cells.Item(27, 8).Value = 6170.3335
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 6170.3335
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 6170.3335
$ws.Cells.Item(27, 14).Value = -6502.3335
$ws.Cells.Item(27, 13).ClearContents()

$ws.Cells.Item(33, 8).Value = 21666.666
$ws.Cells.Item(33, 10).Value = 21666.666
$ws.Cells.Item(33, 12).Value = 21666.666
$ws.Cells.Item(33, 14).Value = -22170.666

$ws.Cells.Item(57, 8).Value = 20000
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 20000
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 20000
$ws.Cells.Item(57, 14).Value = -21640
$ws.Cells.Item(57, 13).ClearContents()

$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).ClearContents()

$ws.Cells.Item(97, 8).Value = 1497.4286
$ws.Cells.Item(97, 9).Value = 421.66666
$ws.Cells.Item(97, 10).Value = 2304.25
$ws.Cells.Item(97, 11).Value = 421.66666
$ws.Cells.Item(97, 12).Value = 2304.25
$ws.Cells.Item(97, 13).Value = 74.33334000000002
$ws.Cells.Item(97, 14).Value = -3296.25

$ws.Cells.Item(126, 8).Value = 6736.875
$ws.Cells.Item(126, 9).Value = 7083
$ws.Cells.Item(126, 10).Value = 5698.5
$ws.Cells.Item(126, 11).Value = 21249
$ws.Cells.Item(126, 12).Value = 17095.5
$ws.Cells.Item(126, 13).Value = -18779
$ws.Cells.Item(126, 14).Value = -22035.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 7125
$ws.Cells.Item(4, 9).Value = 2750
$ws.Cells.Item(4, 10).Value = 11500
$ws.Cells.Item(4, 11).Value = 2750
$ws.Cells.Item(4, 12).Value = 11500
$ws.Cells.Item(4, 13).Value = -2637
$ws.Cells.Item(4, 14).Value = -11726

$ws.Cells.Item(25, 8).Value = 10000
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 13).ClearContents()

$ws.Cells.Item(28, 8).Value = 7125
$ws.Cells.Item(28, 9).Value = 2750
$ws.Cells.Item(28, 10).Value = 11500
$ws.Cells.Item(28, 11).Value = 2750
$ws.Cells.Item(28, 12).Value = 11500
$ws.Cells.Item(28, 13).Value = -2518
$ws.Cells.Item(28, 14).Value = -11964

$ws.Cells.Item(37, 8).Value = 7125
$ws.Cells.Item(37, 9).Value = 2750
$ws.Cells.Item(37, 10).Value = 11500
$ws.Cells.Item(37, 11).Value = 2750
$ws.Cells.Item(37, 12).Value = 11500
$ws.Cells.Item(37, 13).Value = -2643
$ws.Cells.Item(37, 14).Value = -11714

$ws.Cells.Item(40, 8).Value = 5638
$ws.Cells.Item(40, 9).Value = 4866.3335
$ws.Cells.Item(40, 10).Value = 6795.5
$ws.Cells.Item(40, 11).Value = 4866.3335
$ws.Cells.Item(40, 12).Value = 6795.5
$ws.Cells.Item(40, 13).Value = -4730.3335
$ws.Cells.Item(40, 14).Value = -7067.5

$ws.Cells.Item(100, 8).Value = 1000
$ws.Cells.Item(100, 10).Value = 1000
$ws.Cells.Item(100, 12).Value = 1000
$ws.Cells.Item(100, 14).Value = -2082

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1203.0769
$ws.Cells.Item(100, 9).Value = 1164
$ws.Cells.Item(100, 10).Value = 1333.3334
$ws.Cells.Item(100, 11).Value = 2328
$ws.Cells.Item(100, 12).Value = 2666.6668
$ws.Cells.Item(100, 13).Value = -1787
$ws.Cells.Item(100, 14).Value = -3748.6668

$ws.Cells.Item(122, 8).Value = 3300.5
$ws.Cells.Item(122, 9).Value = 1800
$ws.Cells.Item(122, 11).Value = 5400
$ws.Cells.Item(122, 13).Value = -2950
